$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
Write-Host "Last sheet: " $lastSheet.Name
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
Write-Host "Added: " $newSheet.Name
$newSheet.Name = "14_Traditional Knowledge Labels"
Write-Host "Count now: " $wb.Worksheets.Count
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $i ": " $wb.Worksheets.Item($i).Name
}
